# Applies the numeric updates to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 607572.9399999999
$ws.Range("I2").Value = 909503.25
$ws.Range("K2").Value = 909503.25
$ws.Range("M2").Value = -909390.25

$ws.Range("H43").Value = 4094.2307
$ws.Range("I43").Value = 250
$ws.Range("K43").Value = 250
$ws.Range("M43").Value = -181

$ws.Range("H47").Value = 21996
$ws.Range("I47").Value = 21996
$ws.Range("K47").Value = 21996
$ws.Range("M47").Value = -21024

$ws.Range("H62").Value = 14911.385
$ws.Range("I62").Value = 24191.777
$ws.Range("J62").Value = 9998.235000000001
$ws.Range("K62").Value = 24191.777
$ws.Range("L62").Value = 9998.235000000001
$ws.Range("M62").Value = -23567.777
$ws.Range("N62").Value = -11246.235

$ws.Range("H65").Value = 14911.385
$ws.Range("I65").Value = 24191.777
$ws.Range("J65").Value = 9998.235000000001
$ws.Range("K65").Value = 120958.885
$ws.Range("L65").Value = 49991.175
$ws.Range("M65").Value = -117838.885
$ws.Range("N65").Value = -56231.175

$ws.Range("H125").Value = 1547.5
$ws.Range("I125").Value = 750.3333
$ws.Range("K125").Value = 6752.9997
$ws.Range("M125").Value = -4292.9997

$ws.Range("H132").Value = 6065.069
$ws.Range("I132").Value = 4441.391
$ws.Range("K132").Value = 13324.173
$ws.Range("M132").Value = -10794.173

$ws.Range("H137").Value = 6386.7646
$ws.Range("I137").Value = 3309.3
$ws.Range("K137").Value = 9927.900000000001
$ws.Range("M137").Value = -7377.900000000001

$ws.Range("H138").Value = 3906.9404
$ws.Range("I138").Value = 6778.875
$ws.Range("J138").Value = 3231.1912
$ws.Range("K138").Value = 20336.625
$ws.Range("L138").Value = 9693.5736
$ws.Range("M138").Value = -15196.625
$ws.Range("N138").Value = -19973.5736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2645
$ws.Range("I26").Value = 2645
$ws.Range("K26").Value = 2645
$ws.Range("M26").Value = -2315

$ws.Range("H39").Value = 15166.777
$ws.Range("J39").Value = 14500.5
$ws.Range("L39").Value = 14500.5
$ws.Range("N39").Value = -15540.5

$ws.Range("H45").Value = 1433385.2
$ws.Range("I45").Value = 3336000.2
$ws.Range("J45").Value = 6424
$ws.Range("K45").Value = 3336000.2
$ws.Range("L45").Value = 6424
$ws.Range("M45").Value = -3335623.2
$ws.Range("N45").Value = -7178

$ws.Range("H122").Value = 3065
$ws.Range("I122").Value = 2310.9285
$ws.Range("K122").Value = 6932.7855
$ws.Range("M122").Value = -4482.7855

$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 18174
$ws.Range("J21").Value = 18174
$ws.Range("L21").Value = 18174
$ws.Range("N21").Value = -18646

$ws.Range("H23").Value = 2314

$ws.Range("H28").Value = 64000
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H82").Value = 42697
$ws.Range("I82").Value = 12854.6
$ws.Range("K82").Value = 12854.6
$ws.Range("M82").Value = -12471.6

$ws.Range("H85").Value = 42697
$ws.Range("I85").Value = 12854.6
$ws.Range("K85").Value = 12854.6
$ws.Range("M85").Value = -11528.6

$ws.Range("H134").Value = 2643.3635
$ws.Range("I134").Value = 2492.9023
$ws.Range("K134").Value = 7478.706900000001
$ws.Range("M134").Value = -4943.706900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 113.42857
$ws.Range("I7").Value = 109.4
$ws.Range("J7").Value = 115.666664
$ws.Range("K7").Value = 109.4
$ws.Range("L7").Value = 115.666664
$ws.Range("M7").Value = 3.599999999999994
$ws.Range("N7").Value = -341.666664

$ws.Range("H20").Value = 71748.836
$ws.Range("I20").Value = 92000
$ws.Range("J20").Value = 64998.445
$ws.Range("K20").Value = 92000
$ws.Range("L20").Value = 64998.445
$ws.Range("M20").Value = -91764
$ws.Range("N20").Value = -65470.445

$ws.Range("H30").Value = 71748.836
$ws.Range("I30").Value = 92000
$ws.Range("J30").Value = 64998.445
$ws.Range("K30").Value = 92000
$ws.Range("L30").Value = 64998.445
$ws.Range("M30").Value = -91909
$ws.Range("N30").Value = -65180.445

$ws.Range("H31").Value = 5167.1113
$ws.Range("I31").Value = 2726.6667
$ws.Range("J31").Value = 6387.3335
$ws.Range("K31").Value = 2726.6667
$ws.Range("L31").Value = 6387.3335
$ws.Range("M31").Value = -2431.6667
$ws.Range("N31").Value = -6977.3335

$ws.Range("H34").Value = 5167.1113
$ws.Range("I34").Value = 2726.6667
$ws.Range("J34").Value = 6387.3335
$ws.Range("K34").Value = 2726.6667
$ws.Range("L34").Value = 6387.3335
$ws.Range("M34").Value = -2524.6667
$ws.Range("N34").Value = -6791.3335

$ws.Range("H58").Value = 44046.832
$ws.Range("I58").Value = 47832.273
$ws.Range("K58").Value = 47832.273
$ws.Range("M58").Value = -47629.273

$ws.Range("H128").Value = 71748.836
$ws.Range("I128").Value = 92000
$ws.Range("J128").Value = 64998.445
$ws.Range("K128").Value = 92000
$ws.Range("L128").Value = 64998.445
$ws.Range("M128").Value = -87020
$ws.Range("N128").Value = -74958.44500000001

$ws.Range("H134").Value = 38019.715
$ws.Range("I134").Value = 45815.87
$ws.Range("K134").Value = 137447.61
$ws.Range("M134").Value = -134912.61

$ws.Range("H136").Value = 44046.832
$ws.Range("I136").Value = 47832.273
$ws.Range("K136").Value = 143496.819
$ws.Range("M136").Value = -140946.819

$ws.Range("H141").Value = 385962.06
$ws.Range("J141").Value = 410530.78
$ws.Range("L141").Value = 410530.78
$ws.Range("N141").Value = -420890.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 393.4
$ws.Range("I24").Value = 393.4
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1180.2
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H60").Value = 1009.6667
$ws.Range("I60").Value = 1009.6667
$ws.Range("K60").Value = 3029.0001
$ws.Range("M60").Value = -2778.0001

$ws.Range("H107").Value = 532.9
$ws.Range("J107").Value = 531
$ws.Range("L107").Value = 1593
$ws.Range("N107").Value = -5433

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 21500
$ws.Range("J59").Value = 21500
$ws.Range("L59").Value = 21500
$ws.Range("N59").Value = -22666

$ws.Range("H113").Value = 74166.21000000001
$ws.Range("I113").Value = 144809.14
$ws.Range("J113").Value = 3523.2856
$ws.Range("K113").Value = 144809.14
$ws.Range("L113").Value = 3523.2856
$ws.Range("M113").Value = -142639.14
$ws.Range("N113").Value = -7863.2856

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H126").Value = 7142.8276
$ws.Range("I126").Value = 5255.375
$ws.Range("K126").Value = 15766.125
$ws.Range("M126").Value = -13296.125

$ws.Range("H132").Value = 91765.37
$ws.Range("I132").Value = 111945.11
$ws.Range("K132").Value = 335835.33
$ws.Range("M132").Value = -333305.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4274.75
$ws.Range("I40").Value = 3801
$ws.Range("J40").Value = 4559
$ws.Range("K40").Value = 3801
$ws.Range("L40").Value = 4559
$ws.Range("M40").Value = -3665
$ws.Range("N40").Value = -4831

$ws.Range("H55").Value = 1281.3334
$ws.Range("I55").Value = 348.44446
$ws.Range("J55").Value = 2680.6667
$ws.Range("K55").Value = 348.44446
$ws.Range("L55").Value = 2680.6667
$ws.Range("M55").Value = -175.44446
$ws.Range("N55").Value = -3026.6667

$ws.Range("H61").Value = 2961.0476
$ws.Range("J61").Value = 4283.4
$ws.Range("L61").Value = 4283.4
$ws.Range("N61").Value = -4687.4

$ws.Range("H93").Value = 1659.9131
$ws.Range("I93").Value = 1574
$ws.Range("K93").Value = 1574
$ws.Range("M93").Value = -326

$ws.Range("H113").Value = 2961.0476
$ws.Range("J113").Value = 4283.4
$ws.Range("L113").Value = 4283.4
$ws.Range("N113").Value = -8623.4

$ws.Range("H122").Value = 4156.804
$ws.Range("I122").Value = 3563.4814
$ws.Range("K122").Value = 10690.4442
$ws.Range("M122").Value = -8240.4442

$ws.Range("H132").Value = 51170
$ws.Range("I132").Value = 65760.89999999999
$ws.Range("J132").Value = 4965.5
$ws.Range("K132").Value = 197282.7
$ws.Range("L132").Value = 14896.5
$ws.Range("M132").Value = -194752.7
$ws.Range("N132").Value = -19956.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 28790.666
$ws.Range("I54").Value = 21721.111
$ws.Range("K54").Value = 21721.111
$ws.Range("M54").Value = -21201.111

$ws.Range("H81").Value = 773736.4
$ws.Range("I81").Value = 2007.6666
$ws.Range("K81").Value = 4015.3332
$ws.Range("M81").Value = -2954.3332

$ws.Range("H84").Value = 773736.4
$ws.Range("I84").Value = 2007.6666
$ws.Range("K84").Value = 20076.666
$ws.Range("M84").Value = -14772.666

$ws.Range("H126").Value = 47967.047
$ws.Range("I126").Value = 64805.688
$ws.Range("K126").Value = 194417.064
$ws.Range("M126").Value = -191947.064

$ws.Range("H132").Value = 29473.623
$ws.Range("I132").Value = 30583.105
$ws.Range("K132").Value = 91749.315
$ws.Range("M132").Value = -89219.315

